$wb = $excel.ActiveWorkbook

$oldGuid = "adf700cd-e988-4ed0-9a4b-c4ded32d8526"
$newGuid = "81ff3386-6d10-4c04-ad6f-359ca4ef08bd"

$oldHash = "0bbfe7643ee14a706747dac00b15454c0df614d9"
$newHash = "f4f378ba79d76271ad4629a62a122341b2ecc7a1"

# the external hyperlink target address is unchanged by this commit - only
# the cell text / hyperlink display text move to the new guid
$linkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e6b9dcb068198359e26eb6b7107fad8b9acea867/e2e/$oldGuid.md"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-17 10:54:36"

# refresh the hyperlink's display text on B2 (address/relationship id is kept)
$wsOverview.Cells.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $linkAddress, "", "", "e2e\$newGuid.md")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-17 10:54:31"

$wsZhCn.Cells.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $linkAddress, "", "", "$newGuid.md")

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-17 10:54:36"

$wsDeDe.Cells.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $linkAddress, "", "", "$newGuid.md")
